$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the "Upload" row (row 8) for a new "Force" config
# flag, pushing "Upload" and the task-data table down by one row.
$ws.Rows.Item(8).Insert()

# Excel's default row-insert behavior copies the formatting of the row
# above ("Ref", row 7); replicate that, then fill in the new row's values.
$ws.Range("A7:J7").Copy()
$ws.Range("A8:J8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Force"
$ws.Range("B8:J8").Value = $false

# Re-freeze the panes one row lower so the header block (now 10 rows) stays
# pinned above the scrolling task-data table, and leave the selection on
# the (now-shifted) "Upload" row.
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A11").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A9").Select()

$ws.Rows.Item(7).RowHeight = 13.5
$ws.Rows.Item(8).RowHeight = 13.5
$ws.Rows.Item(9).RowHeight = 13.5
